$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("X2").Value = "Utility (Percent)"
$ws.Range("L3").Value = "644 usec"
$ws.Range("M3").Value = "12283k usec"
$ws.Range("N3").Value = "1222478.70 usec"
$ws.Range("O3").Value = "1839 usec"
$ws.Range("P3").Value = "12516k usec"
$ws.Range("Q3").Value = "942546.03 usec"
$ws.Range("L4").Value = "404 usec"
$ws.Range("M4").Value = "5228.2k usec"
$ws.Range("N4").Value = "718579.80 usec"
$ws.Range("O4").Value = "1806 usec"
$ws.Range("P4").Value = "6480.3k usec"
$ws.Range("Q4").Value = "1856615.58 usec"
$ws.Range("L5").Value = "1827 usec"
$ws.Range("M5").Value = "8409.4k usec"
$ws.Range("N5").Value = "1248472.63 usec"
$ws.Range("O5").Value = "1361 usec"
$ws.Range("P5").Value = "8460.8k usec"
$ws.Range("Q5").Value = "103743.50 usec"
$ws.Range("L6").Value = "632 usec"
$ws.Range("M6").Value = "5207.8k usec"
$ws.Range("N6").Value = "719239.21 usec"
$ws.Range("O6").Value = "2 msec"
$ws.Range("P6").Value = "6680 msec"
$ws.Range("Q6").Value = "1870.33 msec"
$ws.Range("L7").Value = "9 msec"
$ws.Range("M7").Value = "3174 msec"
$ws.Range("N7").Value = "1276.35 msec"
$ws.Range("O7").Value = "1775 usec"
$ws.Range("P7").Value = "1696.7k usec"
$ws.Range("Q7").Value = "6417.59 usec"
$ws.Range("L8").Value = "35 msec"
$ws.Range("M8").Value = "2414 msec"
$ws.Range("N8").Value = "1276.09 msec"
$ws.Range("O8").Value = "1615 usec"
$ws.Range("P8").Value = "1079.8k usec"
$ws.Range("Q8").Value = "6875.05 usec"
$ws.Range("L9").Value = "55 msec"
$ws.Range("M9").Value = "2643 msec"
$ws.Range("N9").Value = "1271.05 msec"
$ws.Range("O9").Value = "2 msec"
$ws.Range("P9").Value = "2396 msec"
$ws.Range("Q9").Value = "19.15 msec"
$ws.Range("L10").Value = "634 usec"
$ws.Range("M10").Value = "5010.5k usec"
$ws.Range("N10").Value = "989279.56 usec"
$ws.Range("O10").Value = "2 msec"
$ws.Range("P10").Value = "6141 msec"
$ws.Range("Q10").Value = "971.99 msec"
$ws.Range("L11").Value = "74 msec"
$ws.Range("M11").Value = "2824 msec"
$ws.Range("N11").Value = "1273.74 msec"
$ws.Range("O11").Value = "1907 usec"
$ws.Range("P11").Value = "1807.9k usec"
$ws.Range("Q11").Value = "12720.64 usec"
$ws.Range("L12").Value = "16 msec"
$ws.Range("M12").Value = "2705 msec"
$ws.Range("N12").Value = "1271.44 msec"
$ws.Range("O12").Value = "3 msec"
$ws.Range("P12").Value = "1280 msec"
$ws.Range("Q12").Value = "17.60 msec"
$ws.Range("L13").Value = "346 usec"
$ws.Range("M13").Value = "5199.3k usec"
$ws.Range("N13").Value = "763255.80 usec"
$ws.Range("O13").Value = "1860 usec"
$ws.Range("P13").Value = "6696.6k usec"
$ws.Range("Q13").Value = "1669827.70 usec"
$ws.Range("L14").Value = "10 msec"
$ws.Range("M14").Value = "2693 msec"
$ws.Range("N14").Value = "1272.23 msec"
$ws.Range("O14").Value = "1888 usec"
$ws.Range("P14").Value = "2682.8k usec"
$ws.Range("Q14").Value = "15779.13 usec"
$ws.Range("L15").Value = "18 msec"
$ws.Range("M15").Value = "2617 msec"
$ws.Range("N15").Value = "1220.28 msec"
$ws.Range("O15").Value = "1337 usec"
$ws.Range("P15").Value = "3111.3k usec"
$ws.Range("Q15").Value = "137542.26 usec"
